$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the "K" column (G) values for rows 2-19 to reflect
# the recalculated s_vals (std/mean based K stat, replacing the old Strike# values).
$kValues = @{
    2  = 1
    3  = 4
    4  = 4
    5  = 6
    6  = 6
    7  = 2
    8  = 5
    9  = 3
    10 = 5
    11 = 5
    12 = 7
    13 = 4
    14 = 4
    15 = 3
    16 = 6
    17 = 3
    18 = 2
    19 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
